# Update TC01 (DNBSEQ-G400) queries/filenames to TC07 (Illumina HiSeq 4000)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns D and E hold the same filenames across rows 2-4 (Neo4jData / WebData file names)
$ws.Range("D2:D4").Value = 'TC07_CDS_Filter_InstrumentModel-Illumina HiSeq 4000_Neo4jData.xlsx'
$ws.Range("E2:E4").Value = 'TC07_CDS_Filter_InstrumentModel-Illumina HiSeq 4000_WebData.xlsx'

# Query text column (B) for each tab row, written first so shared strings keep
# the same relative ordering as the authored workbook.
$ws.Range("B2").Value = 'Match (f)<--(g:genomic_info)
WHERE g.instrument_model in [''Illumina HiSeq 4000'']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY `Participant ID`LIMIT 100'
$ws.Range("B3").Value = 'Match (f)<--(g:genomic_info)
WHERE g.instrument_model in [''Illumina HiSeq 4000'']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '''') as `Sample ID`,
 coalesce(p.participant_id,'''') as `Participant ID`,
 coalesce(s.study_name, '''') as `Study Name`,
 coalesce(s.phs_accession,'''') as `Accession`,
coalesce(samp.sample_tumor_status,'''') as `Tumor`,
coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100'
$ws.Range("B4").Value = 'Match (f)<--(g:genomic_info)
WHERE g.instrument_model in [''Illumina HiSeq 4000'']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name, '''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id,'''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER By f.file_name LIMIT 100'

# Stat query (column C) is identical across rows 2-4.
$ws.Range("C2:C4").Value = 'MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in [''Illumina HiSeq 4000'']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files'

# Update selection on the sheet view to D4
$ws.Range("D4").Select()
